$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '40.158.08'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '2.225.65'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue $ws.Range("D5") '293.64'
$ws.Range("E5").Value = '  +1.72%  '
Set-TextValue $ws.Range("D6") '87.94'
$ws.Range("E6").Value = '  +0.52%  '
Set-TextValue $ws.Range("D7") '0.514'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +0.00%  '
Set-TextValue $ws.Range("D9") '0.470'
$ws.Range("E9").Value = '  +0.16%  '
Set-TextValue $ws.Range("D10") '30.70'
$ws.Range("E10").Value = '  +0.60%  '
Set-TextValue $ws.Range("D11") '50.87'
$ws.Range("E11").Value = '  +6.43%  '
Set-TextValue $ws.Range("D12") '0.0783'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("E13").Value = '  +3.39%  '
Set-TextValue $ws.Range("D14") '6.43'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '2.564.21'
$ws.Range("E15").Value = '  +0.25%  '
Set-TextValue $ws.Range("D16") '13.85'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").Value = '2.251.45'
$ws.Range("E17").Value = '  +1.82%  '
Set-TextValue $ws.Range("D18") '0.737'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '40.121.27'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = '0.0₃0890'
$ws.Range("E20").Value = '  +0.73%  '
Set-TextValue $ws.Range("D21") '11.25'
$ws.Range("E21").Value = '  -3.41%  '
Set-TextValue $ws.Range("D22") '5.79'
$ws.Range("E22").Value = '  -0.16%  '
Set-TextValue $ws.Range("D23") '65.67'
$ws.Range("E23").Value = '  +0.33%  '
Set-TextValue $ws.Range("D24") '236.19'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("E27").Value = '  +0.02%  '
Set-TextValue $ws.Range("D28") '23.23'
$ws.Range("E28").Value = '  +2.95%  '
Set-TextValue $ws.Range("D29") '9.34'
$ws.Range("E29").Value = '  +1.36%  '
Set-TextValue $ws.Range("D30") '2.06'
$ws.Range("E30").Value = '  -6.11%  '
Set-TextValue $ws.Range("D31") '158.75'
$ws.Range("E31").Value = '  +3.85%  '
Set-TextValue $ws.Range("D32") '31.92'
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  -0.05%  '
Set-TextValue $ws.Range("D34") '4.97'
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("E35").Value = '  +7.18%  '
Set-TextValue $ws.Range("D36") '0.0715'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("E39").Value = '  +3.82%  '
Set-TextValue $ws.Range("D40") '0.0994'
$ws.Range("E40").Value = '  +0.23%  '
Set-TextValue $ws.Range("D41") '15.69'
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D42").Value = '2.087.97'
$ws.Range("E42").Value = '  -0.67%  '
Set-TextValue $ws.Range("D43") '3.76'
$ws.Range("E43").Value = '  -2.47%  '
Set-TextValue $ws.Range("D44") '19.32'
$ws.Range("E44").Value = '  +9.87%  '
Set-TextValue $ws.Range("D45") '0.0270'
$ws.Range("E45").Value = '  +0.93%  '
Set-TextValue $ws.Range("D46") '10.06'
$ws.Range("E46").Value = '  +1.09%  '
Set-TextValue $ws.Range("D47") '2.75'
$ws.Range("E48").Value = '  -13.31%  '
$ws.Range("D49").Value = '2.423.84'
$ws.Range("E49").Value = '  -0.27%  '
Set-TextValue $ws.Range("D50") '1.47'
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("E51").Value = '  +3.68%  '

Write-Host "Updated cryptos values"
